# Scenario 9 edit: decrease probability of fetal death before 4 weeks
# from conception (Phase1 + Phase2 input tables), and move the active
# worksheet/selection from Phase4 back to Phase1.

$wb = $excel.ActiveWorkbook

# --- Phase2 worksheet: update column C probabilities (E recalculates) ---
$wsPhase2 = $wb.Worksheets.Item("Phase2")
$wsPhase2.Range("C2").Value = 0.1
$wsPhase2.Range("C3").Value = 0.1
$wsPhase2.Range("C4").Value = 0.05
$wsPhase2.Range("C5").Value = 0.05
# Update the saved selection on this sheet (it is not the active tab).
$wsPhase2.Range("C2:C5").Select()

# --- Phase4 worksheet: it loses the "active tab" status (handled below by
# activating Phase1 last), selection stays as saved in the file already.

# --- Phase1 worksheet: update column B probabilities (D recalculates) ---
$wsPhase1 = $wb.Worksheets.Item("Phase1")
$wsPhase1.Range("B2").Value = 0.1
$wsPhase1.Range("B3").Value = 0.1
$wsPhase1.Range("B4").Value = 0.05
$wsPhase1.Range("B5").Value = 0.05

# Make Phase1 the active sheet/tab with the new selection B2:B5.
# Selecting this last ensures Phase1 becomes tabSelected and the
# workbook's active tab, while clearing tabSelected from Phase4.
$wsPhase1.Activate()
$wsPhase1.Range("B2:B5").Select()
